# Updated capital structure database
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D" = -0.0482
    "E" = 0.25
    "G" = 0.1925076656445807
    "H" = 0.1925076656445807
    "I" = 0.1713771497133716
    "J" = 0.1653043610634872
    "K" = 140.7
    "L" = 0.09378749500066656
    "U" = 37.2
    "V" = 0.01829358249323826
    "W" = 0.1811277033985582
    "X" = 0.05517301953945103
    "Y" = 0.1259546838591072
    "Z" = 1.92827763496144
    "AA" = 0.3187527024003129
    "AB" = 0.05487098010626122
    "AC" = 0.2638817222940517
    "AD" = 32.9
    "AF" = 32.9
    "AG" = -4.300000000000004
    "AH" = 0.01592140921409214
    "AI" = 0.03806989122888221
    "AJ" = -0.002119061699191802
    "AK" = -0.005199516324062883
    "AL" = 2.21
    "AM" = 2.21
    "AN" = 0.1243386243386243
    "AO" = 116.3348416289593
    "AP" = -0.01625094482237341
    "AQ" = 116.3348416289593
}

foreach ($col in $updates.Keys) {
    $val = $updates[$col]
    $ws.Range("$col" + "2").Value = $val
    $ws.Range("$col" + "3").Value = $val
}
